# Updates the cryptos list cell values/labels to the latest scraped snapshot.
# Mirrors the upstream "Updated cryptos list ... with GitHub Actions" commit:
# refreshed prices / 1h volume %ages, plus a couple of rows whose coin
# ordering flipped (Bittensor<->PancakeSwap, Mantle<->VeChain).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.824.49"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.674.08"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'600.66"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").Value = "'157.70"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.612"
$ws.Range("E8").Value = "  +3.47%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "'0.400"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "'5.86"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'0.0000202"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "'29.18"
$ws.Range("E14").Value = "  -3.60%  "
$ws.Range("D15").Value = "3.155.09"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "65.661.53"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "2.671.67"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'12.77"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "'4.82"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("D21").Value = "'353.01"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D23").Value = "'69.62"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("D25").Value = "'9.74"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").Value = "'1.68"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("D27").Value = "'1.61"
$ws.Range("E27").Value = "  -3.81%  "
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").Value = "'8.05"
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.14"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'535.11"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "'1.78"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'6.52"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("D35").Value = "'5.49"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("D37").Value = "'20.67"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'158.29"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").Value = "'1.95"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'163.65"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "'4.12"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").Value = "'2.38"
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").Value = "'0.0615"
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("D46").Value = "'22.78"
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0259"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.641"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").Value = "0.0₆0258"
$ws.Range("E49").Value = "  +9.40%  "
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("D51").Value = "'0.0990"
$ws.Range("E51").Value = "  -0.65%  "
